$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 43
$ws.Range("H43").Value = 5059.6665
$ws.Range("I43").Value = 4635.25
$ws.Range("J43").Value = 5214
$ws.Range("K43").Value = 4635.25
$ws.Range("L43").Value = 5214
$ws.Range("M43").Value = -4566.25
$ws.Range("N43").Value = -5352

# Row 82
$ws.Range("H82").Value = 1524.5
$ws.Range("I82").Value = 1524.5
$ws.Range("K82").Value = 4573.5
$ws.Range("M82").Value = -4167.5

# Row 85
$ws.Range("H85").Value = 1524.5
$ws.Range("I85").Value = 1524.5
$ws.Range("K85").Value = 4573.5
$ws.Range("M85").Value = -3169.5

# Row 86
$ws.Range("H86").Value = 14488.385
$ws.Range("I86").Value = 12927.667
$ws.Range("K86").Value = 12927.667
$ws.Range("M86").Value = -11804.667

# Row 89
$ws.Range("H89").Value = 14488.385
$ws.Range("I89").Value = 12927.667
$ws.Range("K89").Value = 64638.335
$ws.Range("M89").Value = -59022.335

# Row 113
$ws.Range("H113").Value = 4997.2
$ws.Range("I113").Value = 3999.75
$ws.Range("J113").Value = 5662.1665
$ws.Range("K113").Value = 3999.75
$ws.Range("L113").Value = 5662.1665
$ws.Range("M113").Value = -745.75
$ws.Range("N113").Value = -12170.1665

# Row 131
$ws.Range("H131").Value = 4409.6
$ws.Range("I131").Value = 3349.6667
$ws.Range("K131").Value = 10049.0001
$ws.Range("M131").Value = -5009.000100000001

# Row 135
$ws.Range("H135").Value = 6808.6
$ws.Range("I135").Value = 1241.3334
$ws.Range("K135").Value = 11172.0006
$ws.Range("M135").Value = -8637.000599999999

# Row 141
$ws.Range("H141").Value = 11429.8
$ws.Range("I141").Value = 6074.5
$ws.Range("K141").Value = 18223.5
$ws.Range("M141").Value = -13043.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("H32").Value = 802065.4
$ws.Range("J32").Value = 26988.5
$ws.Range("L32").Value = 26988.5
$ws.Range("N32").Value = -27562.5

# Row 43
$ws.Range("H43").Value = 74900
$ws.Range("I43").Value = 74900
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 74900
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -74587
$ws.Range("N43").ClearContents()

# Row 45
$ws.Range("H45").Value = 2317.3125
$ws.Range("I45").Value = 1881.4166
$ws.Range("K45").Value = 1881.4166
$ws.Range("M45").Value = -1504.4166

# Row 63
$ws.Range("H63").Value = 2653.3333
$ws.Range("I63").Value = 1974.1
$ws.Range("K63").Value = 1974.1
$ws.Range("M63").Value = -1288.1

# Row 66
$ws.Range("H66").Value = 2653.3333
$ws.Range("I66").Value = 1974.1
$ws.Range("K66").Value = 9870.5
$ws.Range("M66").Value = -6438.5

# Row 74
$ws.Range("H74").Value = 1433141.2
$ws.Range("I74").Value = 2143011
$ws.Range("K74").Value = 2143011
$ws.Range("M74").Value = -2142137

# Row 77
$ws.Range("H77").Value = 1433141.2
$ws.Range("I77").Value = 2143011
$ws.Range("K77").Value = 10715055
$ws.Range("M77").Value = -10710687

# Row 88
$ws.Range("H88").Value = 1976.75
$ws.Range("J88").Value = 1969
$ws.Range("L88").Value = 1969
$ws.Range("N88").Value = -2781

# Row 91
$ws.Range("H91").Value = 1976.75
$ws.Range("J91").Value = 1969
$ws.Range("L91").Value = 1969
$ws.Range("N91").Value = -4777

# Row 132
$ws.Range("H132").Value = 4673.2
$ws.Range("I132").Value = 3506.946
$ws.Range("K132").Value = 10520.838
$ws.Range("M132").Value = -7990.838

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 94
$ws.Range("H94").Value = 3237.303
$ws.Range("I94").Value = 3051.1853
$ws.Range("K94").Value = 3051.1853
$ws.Range("M94").Value = -2600.1853

# Row 134
$ws.Range("H134").Value = 12824260
$ws.Range("I134").Value = 4689.8887
$ws.Range("J134").Value = 41668292
$ws.Range("K134").Value = 14069.6661
$ws.Range("L134").Value = 125004876
$ws.Range("M134").Value = -11534.6661
$ws.Range("N134").Value = -125009946

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 3510550.8
$ws.Range("I31").Value = 3510550.8
$ws.Range("K31").Value = 3510550.8
$ws.Range("M31").Value = -3510255.8

# Row 34
$ws.Range("H34").Value = 3510550.8
$ws.Range("I34").Value = 3510550.8
$ws.Range("K34").Value = 3510550.8
$ws.Range("M34").Value = -3510348.8

# Row 62
$ws.Range("H62").Value = 3826.6667
$ws.Range("I62").Value = 3425
$ws.Range("J62").Value = 4228.3335
$ws.Range("K62").Value = 3425
$ws.Range("L62").Value = 4228.3335
$ws.Range("M62").Value = -2801
$ws.Range("N62").Value = -5476.3335

# Row 65
$ws.Range("H65").Value = 3826.6667
$ws.Range("I65").Value = 3425
$ws.Range("J65").Value = 4228.3335
$ws.Range("K65").Value = 17125
$ws.Range("L65").Value = 21141.6675
$ws.Range("M65").Value = -14005
$ws.Range("N65").Value = -27381.6675

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 102
$ws.Range("H102").Value = 10536.6
$ws.Range("I102").Value = 10254.167
$ws.Range("K102").Value = 10254.167
$ws.Range("M102").Value = -8632.166999999999

# Row 107
$ws.Range("H107").Value = 1637.5
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

# Row 111
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 132
$ws.Range("H132").Value = 14342
$ws.Range("I132").Value = 13792.375
$ws.Range("K132").Value = 41377.125
$ws.Range("M132").Value = -38847.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 3343.32
$ws.Range("J22").Value = 3668.0625
$ws.Range("L22").Value = 3668.0625
$ws.Range("N22").Value = -4258.0625

# Row 27
$ws.Range("H27").Value = 3343.32
$ws.Range("J27").Value = 3668.0625
$ws.Range("L27").Value = 3668.0625
$ws.Range("N27").Value = -3882.0625

# Row 61
$ws.Range("H61").Value = 12376.4
$ws.Range("I61").Value = 13248.75
$ws.Range("K61").Value = 13248.75
$ws.Range("M61").Value = -13046.75

# Row 82
$ws.Range("H82").Value = 1925.45
$ws.Range("I82").Value = 1651.7142
$ws.Range("J82").Value = 2564.1667
$ws.Range("K82").Value = 1651.7142
$ws.Range("L82").Value = 2564.1667
$ws.Range("M82").Value = -1290.7142
$ws.Range("N82").Value = -3286.1667

# Row 85
$ws.Range("H85").Value = 1925.45
$ws.Range("I85").Value = 1651.7142
$ws.Range("J85").Value = 2564.1667
$ws.Range("K85").Value = 1651.7142
$ws.Range("L85").Value = 2564.1667
$ws.Range("M85").Value = -403.7141999999999
$ws.Range("N85").Value = -5060.1667

# Row 88
$ws.Range("H88").Value = 46786.668
$ws.Range("I88").Value = 25171
$ws.Range("J88").Value = 57594.5
$ws.Range("K88").Value = 25171
$ws.Range("L88").Value = 57594.5
$ws.Range("M88").Value = -24743
$ws.Range("N88").Value = -58450.5

# Row 91
$ws.Range("H91").Value = 46786.668
$ws.Range("I91").Value = 25171
$ws.Range("J91").Value = 57594.5
$ws.Range("K91").Value = 25171
$ws.Range("L91").Value = 57594.5
$ws.Range("M91").Value = -23689
$ws.Range("N91").Value = -60558.5

# Row 113
$ws.Range("H113").Value = 12376.4
$ws.Range("I113").Value = 13248.75
$ws.Range("K113").Value = 13248.75
$ws.Range("M113").Value = -11078.75

# Row 122
$ws.Range("H122").Value = 5953.091
$ws.Range("I122").Value = 3496.3333
$ws.Range("K122").Value = 10488.9999
$ws.Range("M122").Value = -8038.999899999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")

# Row 86
$ws.Range("H86").Value = 36000
$ws.Range("J86").Value = 36000
$ws.Range("L86").Value = 36000
$ws.Range("N86").Value = -38246

# Row 89
$ws.Range("H89").Value = 36000
$ws.Range("J89").Value = 36000
$ws.Range("L89").Value = 180000
$ws.Range("N89").Value = -191232

# Row 100
$ws.Range("H100").Value = 1699.6666
$ws.Range("I100").Value = 733
$ws.Range("K100").Value = 1466
$ws.Range("M100").Value = -925

# Row 113
$ws.Range("H113").Value = 3122.3462
$ws.Range("I113").Value = 3108.7856
$ws.Range("J113").Value = 3138.1667
$ws.Range("K113").Value = 9326.356800000001
$ws.Range("L113").Value = 9414.500100000001
$ws.Range("M113").Value = -7156.356800000001
$ws.Range("N113").Value = -13754.5001

# Row 126
$ws.Range("H126").Value = 2966.3333
$ws.Range("I126").Value = 2500
$ws.Range("K126").Value = 7500
$ws.Range("M126").Value = -5030

# Row 132
$ws.Range("H132").Value = 5052374
$ws.Range("I132").Value = 5954112.5
$ws.Range("J132").Value = 2639
$ws.Range("K132").Value = 17862337.5
$ws.Range("L132").Value = 7917
$ws.Range("M132").Value = -17859807.5
$ws.Range("N132").Value = -12977
